# The "Prix Spot" sheet logs one new day-column per update, always inserted
# right before the "01-oct." block (column FA) rather than appended at the
# end. This commit adds "30-dec" as that new column: insert a blank column
# at FA (shifting FA:GE -> FB:GF), write the new header in row 1, and fill
# the as-yet-unpublished data rows (2-25) with "-" placeholders, matching
# the convention used by the other not-yet-available day columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

$ws.Columns("FA:FA").Insert()

$ws.Range("FA1").Value = "30-dec"
$ws.Range("FA2:FA25").Value = "-"
